$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I3").Value = 0.9119513547144661
$ws.Range("J3").Value = 0.2046926104417895
$ws.Range("K3").Value = -0.7920397481549506
$ws.Range("L3").Value = 3.263840307338731
